# Corrections to a couple of mislabeled/misspelled header cells in the
# "Compañias" form layout:
#   D5:  "Código P"  -> "Código postal"
#   B28: "Telefono"  -> "Teléfono"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Código postal"
$ws.Range("B28").Value = "Teléfono"
